$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# ---------------------------------------------------------------------------
# 1. Remove the "Invoice" sheet (its content - the InvoicePage_URL label and
#    the lightning URL - moved into the CreateInvoice sheet; the tab itself
#    is no longer needed in the LUI flow).
# ---------------------------------------------------------------------------
$null = $wb.Worksheets.Item("Invoice").Delete()

# ---------------------------------------------------------------------------
# 2. AddHeader - new sample data (AC5501_CUSTOMER / D1-AC instead of the old
#    DJ-CUST1 customer), widen column A and wrap the long customer label.
# ---------------------------------------------------------------------------
$wsAddHeader = $wb.Worksheets.Item("AddHeader")
$wsAddHeader.Columns.Item(1).ColumnWidth = 41.166666666666664
$wsAddHeader.Range("A2").Value = "AC5501_CUSTOMER (AC5501_CUSTOMER)"
$wsAddHeader.Range("B2").Value = "D1-AC"
$wsAddHeader.Range("A2").WrapText = $true
$wsAddHeader.Rows.Item(2).RowHeight = 19.5
$wsAddHeader.PageSetup.Orientation = 1

# ---------------------------------------------------------------------------
# 3. AddLine - swap the sample Stock/Configured product identifiers for the
#    new Automation5501-1 / Automation5501-2 products.
# ---------------------------------------------------------------------------
$wsAddLine = $wb.Worksheets.Item("AddLine")
$wsAddLine.Range("C2").Value = "Automation5501-1 (Stock-Mfg-LotYesSerialYes)"
$wsAddLine.Range("C3").Value = "Automation5501-2 (Mfg-LotYes)"

# ---------------------------------------------------------------------------
# 4. AddHeader_SOAPI - refreshed record ids / customer for the SOAP API path.
# ---------------------------------------------------------------------------
$wsAddHeaderSOAPI = $wb.Worksheets.Item("AddHeader_SOAPI")
$wsAddHeaderSOAPI.Range("C2").Value = "a5w0W0000019TVB"
$wsAddHeaderSOAPI.Range("D2").Value = "a5w0W0000019TVB"
$wsAddHeaderSOAPI.Range("E2").Value = "a4i0W000001HpNn"
$wsAddHeaderSOAPI.Range("F2").Value = "AC5501_CUSTOMER (AC5501_CUSTOMER)"
$wsAddHeaderSOAPI.Range("G2").Value = "a6Q0W000001Q8bG"

# ---------------------------------------------------------------------------
# 5. AddLine_SOAPI - refreshed record ids / product labels, selection moved.
# ---------------------------------------------------------------------------
$wsAddLineSOAPI = $wb.Worksheets.Item("AddLine_SOAPI")
$wsAddLineSOAPI.Range("C2").Value = "a5w0W0000019TVB"
$wsAddLineSOAPI.Range("D2").Value = "a5w0W0000019TVB"
$wsAddLineSOAPI.Range("E2").Value = "a4i0W000001HpNn"
$wsAddLineSOAPI.Range("F2").Value = "a6Q0W000001Q8bG"
$wsAddLineSOAPI.Range("G2").Value = "a5N6T0000012ymD"
$wsAddLineSOAPI.Range("H2").Value = "Automation5501-1 (Stock-Mfg-LotYesSerialYes)"
$wsAddLineSOAPI.Range("C3").Value = "a5w0W0000019TVB"
$wsAddLineSOAPI.Range("D3").Value = "a5w0W0000019TVB"
$wsAddLineSOAPI.Range("E3").Value = "a4i0W000001HpNn"
$wsAddLineSOAPI.Range("F3").Value = "a6Q0W000001Q8bG"
$wsAddLineSOAPI.Range("G3").Value = "a5N6T0000012ymN"
$wsAddLineSOAPI.Range("H3").Value = "Automation5501-2 (Mfg-LotYes)"

$null = $wsAddLineSOAPI.Activate()
$null = $wsAddLineSOAPI.Range("D8").Select()

# ---------------------------------------------------------------------------
# 6. Leave the workbook with AddHeader_SOAPI as the active / selected tab.
# ---------------------------------------------------------------------------
$null = $wsAddHeaderSOAPI.Activate()
